# B6-PowerPoint.pptx edit script
# 1) Re-style the three summary tables (slides 14, 15, 16) from the deck's
#    custom "Table_0" style to the built-in table style
#    {54301C63-B454-4110-B9A8-908C82A82E35}.
# 2) Re-point the presentation's theme colours from the "Integral" (Red
#    Violet) palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$tableSlides = @(14, 15, 16)
foreach ($slideNum in $tableSlides) {
    $slide = $p.Slides.Item($slideNum)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{54301C63-B454-4110-B9A8-908C82A82E35}")
        }
    }
}

# --- 2) Theme colours --------------------------------------------------
# Office theme palette, in the standard dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink order exposed by ThemeColorScheme.Item(1..12).
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hexColor = $officeColors[$i - 1]
    # PowerPoint COM RGB values are packed as R + G*256 + B*65536.
    $r = ($hexColor -band 0xFF0000) / 0x10000
    $g = ($hexColor -band 0x00FF00) / 0x100
    $b = ($hexColor -band 0x0000FF)
    $comRgb = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($i).RGB = $comRgb
}
